# Add the two missing "Dolar Observado" records (2019-11-20 / 2019-11-21)
# to the valor_dolar sheet, formatting column A as a yyyy-mm-dd date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: 2019-11-20 -> 797.17 ---------------------------------------
$ws.Range("A2").Value = 43789
$ws.Range("A2").Style = "Normal"
$ws.Range("A2").NumberFormat = "yyyy\-mm\-dd"

$ws.Range("B2").Value = 797.17
$ws.Range("B2").Style = "Normal"

$ws.Rows.Item(2).RowHeight = 15

# --- Row 3: 2019-11-21 -> 797.17 ---------------------------------------
$ws.Range("A3").Value = 43790
$ws.Range("A3").Style = "Normal"
$ws.Range("A3").NumberFormat = "yyyy\-mm\-dd"

$ws.Range("B3").Value = 797.17
$ws.Range("B3").Style = "Normal"

$ws.Rows.Item(3).RowHeight = 15

# --- Match the saved selection in the source workbook -------------------
$null = $ws.Range("B7").Select()
